# Insert 3 new price records (rows 759-761) ahead of the existing "Papa"
# data in the Terminal Hortofrutícola Agro Chillán sheet, shifting the
# previously-existing rows 759-796 down to 762-799.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Make room for the new rows by inserting 3 blank rows at row 759.
$ws.Rows("759:761").Insert()

# --- New row 759 ---
$ws.Range("A759").Value = 7
$ws.Range("B759").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C759").Value = "Ñuble"
$ws.Range("D759").Value = 45267
$ws.Range("E759").Value = 16
$ws.Range("F759").Value = 100114001
$ws.Range("G759").Value = "Papa"
$ws.Range("H759").Value = "Asterix"
$ws.Range("I759").Value = "1a nueva(o)"
$ws.Range("J759").Value = 200
$ws.Range("K759").Value = 21000
$ws.Range("L759").Value = 21000
$ws.Range("M759").Value = 21000
$ws.Range("N759").Value = "$/saco 25 kilos"
$ws.Range("O759").Value = "Región del Maule"
$ws.Range("P759").Value = 840
$ws.Range("Q759").Value = 25
$ws.Range("R759").Value = "Hortaliza"

# --- New row 760 ---
$ws.Range("A760").Value = 7
$ws.Range("B760").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C760").Value = "Ñuble"
$ws.Range("D760").Value = 45267
$ws.Range("E760").Value = 16
$ws.Range("F760").Value = 100114001
$ws.Range("G760").Value = "Papa"
$ws.Range("H760").Value = "Asterix"
$ws.Range("I760").Value = "2a nueva(o)"
$ws.Range("J760").Value = 200
$ws.Range("K760").Value = 19000
$ws.Range("L760").Value = 19000
$ws.Range("M760").Value = 19000
$ws.Range("N760").Value = "$/saco 25 kilos"
$ws.Range("O760").Value = "Región del Maule"
$ws.Range("P760").Value = 760
$ws.Range("Q760").Value = 25
$ws.Range("R760").Value = "Hortaliza"

# --- New row 761 ---
$ws.Range("A761").Value = 7
$ws.Range("B761").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C761").Value = "Ñuble"
$ws.Range("D761").Value = 45267
$ws.Range("E761").Value = 16
$ws.Range("F761").Value = 100114001
$ws.Range("G761").Value = "Papa"
$ws.Range("H761").Value = "Rosara"
$ws.Range("I761").Value = "1a nueva(o)"
$ws.Range("J761").Value = 400
$ws.Range("K761").Value = 19000
$ws.Range("L761").Value = 20000
$ws.Range("M761").Value = 19500
$ws.Range("N761").Value = "$/saco 25 kilos"
$ws.Range("O761").Value = "Región del Maule"
$ws.Range("P761").Value = 780
$ws.Range("Q761").Value = 25
$ws.Range("R761").Value = "Hortaliza"

# Apply the date number format used by the rest of column D (Fecha) to
# the freshly-inserted cells so the values round-trip as dates.
$ws.Range("D759:D761").NumberFormat = "YYYY-MM-DD HH:MM:SS"
